$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change word English to Spanish / add missing accents (label fixes)
$ws.Range("D1").Value = "Correo electrónico (*)"
$ws.Range("F1").Value = "Subárea 1 (*)"
$ws.Range("G1").Value = "Subárea 2"
$ws.Range("H1").Value = "Subárea 3"

# Preserve gridlines visibility (engine defaults to hidden otherwise)
$excel.ActiveWindow.DisplayGridlines = $true

# hotfix in view list users: move selection/scroll to Q1
$ws.Range("Q1").Select()
$excel.ActiveWindow.ScrollColumn = 10
$excel.ActiveWindow.ScrollRow = 1
